$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "2025-08-05 12:19:32"
$ws.Range("B12").Value = "create-repo"
$ws.Range("C12").Value = "new-organization97"
$ws.Range("H12").Value = "deerepo"

# Plain Value = "False" is auto-typed to a Boolean by the engine (mirrors
# Excel's "typed into a cell" behavior). The source workbook stores this
# column as literal text, so route the text through a formula + paste-special
# (values only) to land a literal text "False" string in the cell instead of
# a Boolean.
$ws.Range("Z1").Formula = "=""False"""
$ws.Range("Z1").Copy()
$ws.Range("I12").PasteSpecial(-4163)
$ws.Range("Z1").Value = ""
